$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3947
$ws.Range("J51").Value = 3947
$ws.Range("L51").Value = 3947
$ws.Range("N51").Value = -4915
$ws.Range("H95").Value = 17500
$ws.Range("J95").Value = 17500
$ws.Range("L95").Value = 17500
$ws.Range("N95").Value = -22992
$ws.Range("H98").Value = 1376.4117
$ws.Range("I98").Value = 1376.4117
$ws.Range("K98").Value = 1376.4117
$ws.Range("M98").Value = 121.5882999999999
$ws.Range("H113").Value = 2282.2727
$ws.Range("I113").Value = 2166.6667
$ws.Range("J113").Value = 2300.5264
$ws.Range("K113").Value = 2166.6667
$ws.Range("L113").Value = 2300.5264
$ws.Range("M113").Value = 1087.3333
$ws.Range("N113").Value = -8808.526400000001
$ws.Range("H122").Value = 1376.4117
$ws.Range("I122").Value = 1376.4117
$ws.Range("K122").Value = 4129.2351
$ws.Range("M122").Value = -1679.2351

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 14199.5
$ws.Range("J24").Value = 14199.5
$ws.Range("L24").Value = 14199.5
$ws.Range("N24").Value = -14947.5
$ws.Range("H61").Value = 1363.7
$ws.Range("I61").Value = 1232.9412
$ws.Range("K61").Value = 1232.9412
$ws.Range("M61").Value = -1020.9412
$ws.Range("H74").Value = 631.2292
$ws.Range("I74").Value = 575.8461
$ws.Range("K74").Value = 575.8461
$ws.Range("M74").Value = 298.1539
$ws.Range("H77").Value = 631.2292
$ws.Range("I77").Value = 575.8461
$ws.Range("K77").Value = 2879.2305
$ws.Range("M77").Value = 1488.7695
$ws.Range("H100").Value = 14199.5
$ws.Range("J100").Value = 14199.5
$ws.Range("L100").Value = 14199.5
$ws.Range("N100").Value = -16363.5
$ws.Range("H136").Value = 1363.7
$ws.Range("I136").Value = 1232.9412
$ws.Range("K136").Value = 3698.8236
$ws.Range("M136").Value = -1148.8236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 29943.5
$ws.Range("J95").Value = 29943.5
$ws.Range("L95").Value = 29943.5
$ws.Range("N95").Value = -35435.5
$ws.Range("H106").Value = 19249.5
$ws.Range("J106").Value = 19249.5
$ws.Range("L106").Value = 19249.5
$ws.Range("N106").Value = -21773.5
$ws.Range("H126").Value = 49966.668
$ws.Range("J126").Value = 49966.668
$ws.Range("L126").Value = 49966.668
$ws.Range("N126").Value = -59846.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 1234
$ws.Range("I44").Value = 1234
$ws.Range("K44").Value = 1234
$ws.Range("M44").Value = -792
$ws.Range("H55").Value = 3000
$ws.Range("I55").Value = 3000
$ws.Range("K55").Value = 3000
$ws.Range("M55").Value = -2685
$ws.Range("H58").Value = 884.875
$ws.Range("I58").Value = 854.2941
$ws.Range("J58").Value = 959.1429000000001
$ws.Range("K58").Value = 854.2941
$ws.Range("L58").Value = 959.1429000000001
$ws.Range("M58").Value = -651.2941
$ws.Range("N58").Value = -1365.1429
$ws.Range("H115").Value = 47499
$ws.Range("J115").Value = 47499
$ws.Range("L115").Value = 47499
$ws.Range("N115").Value = -49849
$ws.Range("H136").Value = 884.875
$ws.Range("I136").Value = 854.2941
$ws.Range("J136").Value = 959.1429000000001
$ws.Range("K136").Value = 2562.8823
$ws.Range("L136").Value = 2877.4287
$ws.Range("M136").Value = -12.88229999999976
$ws.Range("N136").Value = -7977.4287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 422546.3
$ws.Range("I4").Value = 62639.312
$ws.Range("J4").Value = 782453.3
$ws.Range("K4").Value = 187917.936
$ws.Range("L4").Value = 2347359.9
$ws.Range("M4").Value = -187805.936
$ws.Range("N4").Value = -2347583.9
$ws.Range("H68").Value = 1474.0435
$ws.Range("I68").Value = 651
$ws.Range("K68").Value = 1953
$ws.Range("M68").Value = -1142
$ws.Range("H71").Value = 1474.0435
$ws.Range("I71").Value = 651
$ws.Range("K71").Value = 5859
$ws.Range("M71").Value = -1803
$ws.Range("H131").Value = 37038690
$ws.Range("I131").Value = 111111770
$ws.Range("J131").Value = 2150.111
$ws.Range("K131").Value = 333335310
$ws.Range("L131").Value = 6450.333
$ws.Range("M131").Value = -333330270
$ws.Range("N131").Value = -16530.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 22999.5
$ws.Range("J98").Value = 22999.5
$ws.Range("L98").Value = 22999.5
$ws.Range("N98").Value = -28989.5
$ws.Range("H102").Value = 1635.6842
$ws.Range("I102").Value = 1644.1538
$ws.Range("J102").Value = 1617.3334
$ws.Range("K102").Value = 1644.1538
$ws.Range("L102").Value = 1617.3334
$ws.Range("M102").Value = -22.15380000000005
$ws.Range("N102").Value = -4861.3334
$ws.Range("H122").Value = 2750.3845
$ws.Range("I122").Value = 1622.4286
$ws.Range("K122").Value = 4867.2858
$ws.Range("M122").Value = -2417.2858
$ws.Range("H126").Value = 2525.6924
$ws.Range("I126").Value = 1704.8572
$ws.Range("J126").Value = 3483.3333
$ws.Range("K126").Value = 5114.571599999999
$ws.Range("L126").Value = 10449.9999
$ws.Range("M126").Value = -2644.571599999999
$ws.Range("N126").Value = -15389.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1263.8235
$ws.Range("I22").Value = 1165
$ws.Range("J22").Value = 1405
$ws.Range("K22").Value = 1165
$ws.Range("L22").Value = 1405
$ws.Range("M22").Value = -870
$ws.Range("N22").Value = -1995
$ws.Range("H27").Value = 1263.8235
$ws.Range("I27").Value = 1165
$ws.Range("J27").Value = 1405
$ws.Range("K27").Value = 1165
$ws.Range("L27").Value = 1405
$ws.Range("M27").Value = -1058
$ws.Range("N27").Value = -1619
$ws.Range("H46").Value = 1900
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("H55").Value = 651.8570999999999
$ws.Range("I55").Value = 129.6
$ws.Range("J55").Value = 942
$ws.Range("K55").Value = 129.6
$ws.Range("L55").Value = 942
$ws.Range("M55").Value = 43.40000000000001
$ws.Range("N55").Value = -1288
$ws.Range("H82").Value = 2430.8333
$ws.Range("I82").Value = 2660
$ws.Range("J82").Value = 2201.6667
$ws.Range("K82").Value = 2660
$ws.Range("L82").Value = 2201.6667
$ws.Range("M82").Value = -2299
$ws.Range("N82").Value = -2923.6667
$ws.Range("H85").Value = 2430.8333
$ws.Range("I85").Value = 2660
$ws.Range("J85").Value = 2201.6667
$ws.Range("K85").Value = 2660
$ws.Range("L85").Value = 2201.6667
$ws.Range("M85").Value = -1412
$ws.Range("N85").Value = -4697.6667
$ws.Range("H97").Value = 15000
$ws.Range("J97").Value = 15000
$ws.Range("L97").Value = 15000
$ws.Range("N97").Value = -16982
$ws.Range("H122").Value = 19232982
$ws.Range("I122").Value = 35716310
$ws.Range("J122").Value = 2430.8333
$ws.Range("K122").Value = 107148930
$ws.Range("L122").Value = 7292.499899999999
$ws.Range("M122").Value = -107146480
$ws.Range("N122").Value = -12192.4999
$ws.Range("H132").Value = 65099.5
$ws.Range("I132").Value = 2218.9
$ws.Range("K132").Value = 6656.700000000001
$ws.Range("M132").Value = -4126.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H132").Value = 2737.111
$ws.Range("I132").Value = 2271.3333
$ws.Range("K132").Value = 6813.999899999999
$ws.Range("M132").Value = -4283.999899999999
